$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns (new reporting quarters) before column D, shifting existing
# quarterly data from D:K to F:M.
$ws.Range("D1:E1").EntireColumn.Insert()

# The newly inserted columns inherit formatting from column C by default; copy the
# number formatting from column F (the old column D, now shifted right) onto the new
# D:E columns so date/number styles match the rest of the table.
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new quarters' worth of data in columns D and E.
$ws.Cells.Item(7,4).Value = 43465
$ws.Cells.Item(7,5).Value = 43373
$ws.Cells.Item(8,4).Value = 74200
$ws.Cells.Item(8,5).Value = 71700
$ws.Cells.Item(9,4).Value = 17500
$ws.Cells.Item(9,5).Value = 16500
$ws.Cells.Item(10,4).Value = 56700
$ws.Cells.Item(10,5).Value = 55200
$ws.Cells.Item(12,4).Value = 15100
$ws.Cells.Item(12,5).Value = 12500
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 0
$ws.Cells.Item(14,4).Value = 0
$ws.Cells.Item(14,5).Value = 0
$ws.Cells.Item(15,4).Value = 0
$ws.Cells.Item(15,5).Value = 0
$ws.Cells.Item(17,4).Value = 61300
$ws.Cells.Item(17,5).Value = 53500
$ws.Cells.Item(18,4).Value = 12900
$ws.Cells.Item(18,5).Value = 18200
$ws.Cells.Item(20,4).Value = 1900
$ws.Cells.Item(20,5).Value = 1100
$ws.Cells.Item(21,4).Value = 22500
$ws.Cells.Item(21,5).Value = 26200
$ws.Cells.Item(22,4).Value = 100
$ws.Cells.Item(22,5).Value = 0
$ws.Cells.Item(23,4).Value = 14800
$ws.Cells.Item(23,5).Value = 19200
$ws.Cells.Item(24,4).Value = 400
$ws.Cells.Item(24,5).Value = -4200
$ws.Cells.Item(25,4).Value = 0
$ws.Cells.Item(25,5).Value = 0
$ws.Cells.Item(26,4).Value = 14400
$ws.Cells.Item(26,5).Value = 23500
$ws.Cells.Item(27,4).Value = 14400
$ws.Cells.Item(27,5).Value = 23500
$ws.Cells.Item(28,4).Value = 0
$ws.Cells.Item(28,5).Value = 0
$ws.Cells.Item(29,4).Value = "NA"
$ws.Cells.Item(29,5).Value = "NA"
$ws.Cells.Item(30,4).Value = 0
$ws.Cells.Item(30,5).Value = 0
$ws.Cells.Item(31,4).Value = 0
$ws.Cells.Item(31,5).Value = 0
$ws.Cells.Item(32,4).Value = -1900
$ws.Cells.Item(32,5).Value = -1100
$ws.Cells.Item(33,4).Value = 14400
$ws.Cells.Item(33,5).Value = 23500
$ws.Cells.Item(34,4).Value = 0
$ws.Cells.Item(34,5).Value = 0
$ws.Cells.Item(35,4).Value = 14400
$ws.Cells.Item(35,5).Value = 23500
$ws.Cells.Item(38,4).Value = 43465
$ws.Cells.Item(38,5).Value = 43373
$ws.Cells.Item(41,4).Value = 41000
$ws.Cells.Item(41,5).Value = 98000
$ws.Cells.Item(42,4).Value = 248100
$ws.Cells.Item(42,5).Value = 228500
$ws.Cells.Item(43,4).Value = 75800
$ws.Cells.Item(43,5).Value = 58600
$ws.Cells.Item(44,4).Value = 0
$ws.Cells.Item(44,5).Value = 0
$ws.Cells.Item(45,4).Value = 14000
$ws.Cells.Item(45,5).Value = 16600
$ws.Cells.Item(46,4).Value = 379000
$ws.Cells.Item(46,5).Value = 401700
$ws.Cells.Item(47,4).Value = 79200
$ws.Cells.Item(47,5).Value = 66600
$ws.Cells.Item(48,4).Value = 61400
$ws.Cells.Item(48,5).Value = 64500
$ws.Cells.Item(49,4).Value = 29200
$ws.Cells.Item(49,5).Value = 15600
$ws.Cells.Item(50,4).Value = 0
$ws.Cells.Item(50,5).Value = 0
$ws.Cells.Item(51,4).Value = 0
$ws.Cells.Item(51,5).Value = 0
$ws.Cells.Item(52,4).Value = 36900
$ws.Cells.Item(52,5).Value = 37600
$ws.Cells.Item(53,4).Value = 0
$ws.Cells.Item(53,5).Value = 0
$ws.Cells.Item(54,4).Value = 585700
$ws.Cells.Item(54,5).Value = 586000
$ws.Cells.Item(57,4).Value = 5600
$ws.Cells.Item(57,5).Value = 4700
$ws.Cells.Item(58,4).Value = 1600
$ws.Cells.Item(58,5).Value = "NA"
$ws.Cells.Item(59,4).Value = 189800
$ws.Cells.Item(59,5).Value = 178300
$ws.Cells.Item(60,4).Value = 196900
$ws.Cells.Item(60,5).Value = 183000
$ws.Cells.Item(61,4).Value = 0
$ws.Cells.Item(61,5).Value = 0
$ws.Cells.Item(62,4).Value = 30800
$ws.Cells.Item(62,5).Value = 28800
$ws.Cells.Item(63,4).Value = 0
$ws.Cells.Item(63,5).Value = 0
$ws.Cells.Item(64,4).Value = 0
$ws.Cells.Item(64,5).Value = 0
$ws.Cells.Item(65,4).Value = 0
$ws.Cells.Item(65,5).Value = 0
$ws.Cells.Item(66,4).Value = 227700
$ws.Cells.Item(66,5).Value = 211900
$ws.Cells.Item(68,4).Value = 0
$ws.Cells.Item(68,5).Value = 0
$ws.Cells.Item(69,4).Value = 0
$ws.Cells.Item(69,5).Value = 0
$ws.Cells.Item(70,4).Value = 0
$ws.Cells.Item(70,5).Value = 0
$ws.Cells.Item(71,4).Value = 0
$ws.Cells.Item(71,5).Value = 0
$ws.Cells.Item(72,4).Value = 28000
$ws.Cells.Item(72,5).Value = 45800
$ws.Cells.Item(73,4).Value = 0
$ws.Cells.Item(73,5).Value = 0
$ws.Cells.Item(74,4).Value = 0
$ws.Cells.Item(74,5).Value = 0
$ws.Cells.Item(75,4).Value = 0
$ws.Cells.Item(75,5).Value = 0
$ws.Cells.Item(76,4).Value = 358000
$ws.Cells.Item(76,5).Value = 374100
$ws.Cells.Item(77,4).Value = 0
$ws.Cells.Item(77,5).Value = 0
$ws.Cells.Item(80,4).Value = 43465
$ws.Cells.Item(80,5).Value = 43373
$ws.Cells.Item(81,4).Value = 14400
$ws.Cells.Item(81,5).Value = 23500
$ws.Cells.Item(83,4).Value = 7700
$ws.Cells.Item(83,5).Value = 7000
$ws.Cells.Item(84,4).Value = 0
$ws.Cells.Item(84,5).Value = 0
$ws.Cells.Item(85,4).Value = 0
$ws.Cells.Item(85,5).Value = 0
$ws.Cells.Item(86,4).Value = 0
$ws.Cells.Item(86,5).Value = 0
$ws.Cells.Item(87,4).Value = 0
$ws.Cells.Item(87,5).Value = 0
$ws.Cells.Item(88,4).Value = 0
$ws.Cells.Item(88,5).Value = 0
$ws.Cells.Item(89,4).Value = 26600
$ws.Cells.Item(89,5).Value = 31600
$ws.Cells.Item(91,4).Value = -3300
$ws.Cells.Item(91,5).Value = -6300
$ws.Cells.Item(92,4).Value = 0
$ws.Cells.Item(92,5).Value = 0
$ws.Cells.Item(93,4).Value = 0
$ws.Cells.Item(93,5).Value = 0
$ws.Cells.Item(94,4).Value = -45000
$ws.Cells.Item(94,5).Value = 1500
$ws.Cells.Item(96,4).Value = 0
$ws.Cells.Item(96,5).Value = 0
$ws.Cells.Item(97,4).Value = 0
$ws.Cells.Item(97,5).Value = 0
$ws.Cells.Item(98,4).Value = 0
$ws.Cells.Item(98,5).Value = 0
$ws.Cells.Item(99,4).Value = 0
$ws.Cells.Item(99,5).Value = 0
$ws.Cells.Item(100,4).Value = -38600
$ws.Cells.Item(100,5).Value = -21900
$ws.Cells.Item(101,4).Value = 0
$ws.Cells.Item(101,5).Value = 0
$ws.Cells.Item(102,4).Value = -56900
$ws.Cells.Item(102,5).Value = 11200

# Row 58 ("Short/Current Long Term Debt") source refresh also marked several of the
# older quarters as not-available.
$ws.Cells.Item(58,6).Value = "NA"
$ws.Cells.Item(58,7).Value = "NA"
$ws.Cells.Item(58,8).Value = "NA"
$ws.Cells.Item(58,9).Value = "NA"
$ws.Cells.Item(58,10).Value = "NA"
